{"js": "// Split the \"Programa\" paragraphs (Portuguese and English) so that the\n// numbered items (\"1. ...\" / \"2. ...\") each start on a new line, using\n// manual line breaks (<w:br/>) instead of running on as one block of text.\n// We do this by finding the exact seams in the text and replacing them\n// with the same text plus a vertical-tab (\\u000b), which Word's text model\n// represents as a line break.\n\nconst body = context.document.body;\n\n// --- Portuguese paragraph -------------------------------------------------\nconst ptBreak1 = body.search(\"contextos: 1.\", { matchCase: true });\nawait context.sync();\nif (ptBreak1.items.length > 0) {\n  ptBreak1.items[0].insertText(\"contextos: \\u000b1.\", \"Replace\");\n  await context.sync();\n}\n\nconst ptBreak2 = body.search(\"Engenharia; 2.\", { matchCase: true });\nawait context.sync();\nif (ptBreak2.items.length > 0) {\n  ptBreak2.items[0].insertText(\"Engenharia; \\u000b2.\", \"Replace\");\n  await context.sync();\n}\n\n// --- English (italic) paragraph ------------------------------------------\nconst enBreak1 = body.search(\"contexts:1.\", { matchCase: true });\nawait context.sync();\nif (enBreak1.items.length > 0) {\n  enBreak1.items[0].insertText(\"contexts:\\u000b1.\", \"Replace\");\n  await context.sync();\n}\n\nconst enBreak2 = body.search(\"course;2.\", { matchCase: true });\nawait context.sync();\nif (enBreak2.items.length > 0) {\n  enBreak2.items[0].insertText(\"course;\\u000b2.\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Split the \"Programa\" paragraphs (Portuguese and English) so that the\n# numbered items (\"1. ...\" / \"2. ...\") each start on a new line, using a\n# manual line break (Word's \"^l\" / <w:br/>) instead of running on as one\n# block of text.\n\n$d = $word.ActiveDocument\n\nfunction Insert-LineBreak($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# Portuguese paragraph (Programa)\nInsert-LineBreak \"contextos: 1.\" \"contextos: ^l1.\"\nInsert-LineBreak \"Engenharia; 2.\" \"Engenharia; ^l2.\"\n\n# English (italic) paragraph\nInsert-LineBreak \"contexts:1.\" \"contexts:^l1.\"\nInsert-LineBreak \"course;2.\" \"course;^l2.\"\n"}
